$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.615.17"
$ws.Range("E2").Value = "  +3.03%  "
$ws.Range("D3").Value = "2.285.35"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'250.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'0.633"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D7").Value = "'72.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.49%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.643"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.18%  "
$ws.Range("D10").Value = "'38.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").Value = "'0.0967"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").Value = "'59.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "'7.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.94%  "
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").Value = "2.628.22"
$ws.Range("E15").Value = "  +4.39%  "
$ws.Range("D16").Value = "'15.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").Value = "'0.883"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.29%  "
$ws.Range("D18").Value = "2.278.64"
$ws.Range("E18").Value = "  +4.07%  "
$ws.Range("D19").Value = "42.618.03"
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").Value = "0.0₃0997"
$ws.Range("E20").Value = "  +4.95%  "
$ws.Range("D21").Value = "'6.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.04%  "
$ws.Range("D22").Value = "'72.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").Value = "'2.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.15%  "
$ws.Range("D24").Value = "'235.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "'3.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("D26").Value = "'11.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").Value = "'167.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "'21.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.23%  "
$ws.Range("D33").Value = "'6.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.42%  "
$ws.Range("D34").Value = "'0.128"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.73%  "
$ws.Range("D35").Value = "'31.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +22.95%  "
$ws.Range("D36").Value = "'0.0802"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("D38").Value = "'4.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.30%  "
$ws.Range("D39").Value = "'4.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.94%  "
$ws.Range("D40").Value = "'0.0307"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "'13.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.66%  "
$ws.Range("E42").Value = "  +6.45%  "
$ws.Range("D43").Value = "'6.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.52%  "
$ws.Range("D44").Value = "'0.212"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.01%  "
$ws.Range("D45").Value = "'9.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.11%  "
$ws.Range("D46").Value = "'62.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").Value = "'4.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("D48").Value = "'0.103"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'1.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("D51").Value = "'96.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.41%  "
